# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 02:21"

# Updated case counts (rows keep their sort order except Guinea-Bisau and
# Eslovenia, which swap positions because Guinea-Bisau's totals overtake
# Eslovenia's).

# Row 4 - Rusia
$ws.Range("B4").Value = 3291304
$ws.Range("C4").Value = 71305
$ws.Range("D4").Value = 1454285
$ws.Range("E4").Value = 1700372
$ws.Range("G4").Value = 825
$ws.Range("H4").Value = 136647

# Row 5 - Peru
$ws.Range("B5").Value = 1804338
$ws.Range("C5").Value = 45235
$ws.Range("E5").Value = 548218
$ws.Range("G5").Value = 1270
$ws.Range("H5").Value = 70524

# Row 6 - Chile
$ws.Range("B6").Value = 822603
$ws.Range("C6").Value = 27761
$ws.Range("E6").Value = 284253

# Row 23 - Canada
$ws.Range("B23").Value = 107126
$ws.Range("C23").Value = 321
$ws.Range("E23").Value = 27466

# Row 77
$ws.Range("B77").Value = 8974
$ws.Range("C77").Value = 9
$ws.Range("E77").Value = 584

# Row 122 now becomes Guinea-Bisau with updated totals (previously Eslovenia)
$ws.Range("A122").Value = "Guinea-Bisau"
$ws.Range("B122").Value = 1842
$ws.Range("C122").Value = 52
$ws.Range("D122").Value = 773
$ws.Range("E122").Value = 1043
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 26

# Row 123 now becomes Eslovenia (previously Guinea-Bisau), carrying the
# values Eslovenia had before this update
$ws.Range("A123").Value = "Eslovenia"
$ws.Range("B123").Value = 1793
$ws.Range("C123").Value = 17
$ws.Range("D123").Value = 1429
$ws.Range("E123").Value = 253
$ws.Range("H123").Value = 111

# Row 161
$ws.Range("D161").Value = 350
$ws.Range("E161").Value = 19

# Row 177
$ws.Range("B177").Value = 150
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 4

# Row 219
$ws.Range("B219").Value = 2
$ws.Range("C219").Value = 1
$ws.Range("E219").Value = 1
